# Add an auto-advance slide transition (<p:transition advTm="…"/>) to every
# slide in the deck. PowerPoint's object model exposes this via
# SlideShowTransition.AdvanceOnTime / .AdvanceTime, where AdvanceTime is
# expressed in seconds (the OOXML advTm attribute is milliseconds).

$p = $ppt.ActivePresentation

# Advance time per slide index, in milliseconds (matches each slide's
# p:transition/@advTm in the target deck). Most slides use the default
# 3000 ms (3 s); a few have bespoke values.
$advanceTimesMs = @{
    1  = 2024
    2  = 1069
    3  = 88984
    4  = 3000
    5  = 3000
    6  = 3000
    7  = 3000
    8  = 3000
    9  = 3000
    10 = 3000
    11 = 3000
    12 = 3000
    13 = 3000
    14 = 3000
    15 = 3000
    16 = 3000
    17 = 3000
    18 = 3000
    19 = 3000
    20 = 3000
    21 = 3000
    22 = 3000
    23 = 3000
    24 = 3000
    25 = 3000
    26 = 3000
    27 = 3000
    28 = 3000
    29 = 3000
    30 = 3000
    31 = 3000
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $ms = $advanceTimesMs[$i]
    if ($null -eq $ms) {
        $ms = 3000
    }
    $slide.SlideShowTransition.AdvanceOnTime = $true
    $slide.SlideShowTransition.AdvanceTime = $ms / 1000.0
}

Write-Output "Applied slide transition advance times to $($p.Slides.Count) slides"
